$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Effort R 1.0")

# New row 34
$ws.Range("A34").Value = 41471
$ws.Range("A34").NumberFormat = "ddd\ dd/mm/yyyy"
$ws.Range("B34").Value = 1
$ws.Range("D34").Value = "Setup creation, prerelease sent to Sudar Muthu"

# New row 35
$ws.Range("A35").Value = 41472
$ws.Range("A35").NumberFormat = "ddd\ dd/mm/yyyy"
$ws.Range("B35").Value = 1.25
$ws.Range("D35").Value = "Revision of manual"

$ws.Range("B35").Select()
